# Populate the "Neutered Status" (column I) values that were missing from
# rows 86-145 of the WebDataCanine sheet. Every row in that range gets "Yes"
# except row 104 (case NCATS-COP01CCB020031), which gets "No".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 86; $r -le 145; $r++) {
    if ($r -eq 104) {
        $ws.Cells.Item($r, 9).Value = "No"
    } else {
        $ws.Cells.Item($r, 9).Value = "Yes"
    }
}
